# Daily attendance processing - reorder "Recorded By" entries in column G.
# Rule observed in the target diff: for each cell in column G (data rows),
# split the comma-separated list of recorders. If there are exactly two
# entries, swap them. If there are three (or more) entries, keep the first
# entry fixed in place and reverse the remaining entries. Cells with a
# single entry are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }
    if ($val -eq "") { continue }

    $parts = $val -split ",\s*"
    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    if ($trimmed.Count -eq 2) {
        $newParts = @($trimmed[1], $trimmed[0])
        $cell.Value2 = [string]::Join(", ", $newParts)
    }
    elseif ($trimmed.Count -ge 3) {
        # keep first entry fixed, reverse the remaining entries (manual
        # reverse loop -- [array]::Reverse() does not mutate reliably here)
        $revRest = @()
        for ($i = $trimmed.Count - 1; $i -ge 1; $i--) {
            $revRest += $trimmed[$i]
        }
        $newParts = @($trimmed[0]) + $revRest
        $cell.Value2 = [string]::Join(", ", $newParts)
    }
    # single-entry (or empty) cells are left untouched
}
